# Update cryptocurrency price/volume figures per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "54.383.08"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.281.99"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +2.58%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.14%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "499.11"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +1.67%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "128.25"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +1.27%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.26%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.52%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0958"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +2.85%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +1.69%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +3.61%  "
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +1.69%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.687.59"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +2.47%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "22.50"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +4.80%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "54.315.07"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +0.58%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +0.76%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.279.35"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +1.53%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.28"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +5.41%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +2.39%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "304.99"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +2.46%  "
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +3.18%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "62.08"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -2.74%  "
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.388.82"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +2.62%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +3.24%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +3.20%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "173.61"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +3.03%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0688"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +2.17%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "17.79"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +2.23%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.927"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +10.47%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +1.15%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +4.57%  "
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +2.20%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "126.03"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -1.77%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +0.91%  "
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +2.86%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +0.93%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "239.71"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -1.17%  "
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +1.82%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "16.33"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +1.23%  "
$c.Style = "Normal"
